$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 531, shifting existing rows 531-592 down to 532-593
# (mirrors Excel's Rows("531:531").Insert Shift:=xlShiftDown)
$ws.Rows.Item(531).Insert(-4121)

# Populate the newly inserted row 531 with the new weekly price observation.
# Columns that keep the same value as the rest of this market/variety block
# (A, B, C, E, F, G, H, I, N, O, Q, R) are re-asserted explicitly too, so the
# row is fully defined regardless of what Insert happened to carry down.
$ws.Cells.Item(531, 1).Value = 10
$ws.Cells.Item(531, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(531, 3).Value = "La Araucanía"
$ws.Cells.Item(531, 4).Value = 45077
$ws.Cells.Item(531, 5).Value = 9
$ws.Cells.Item(531, 6).Value = 100114014
$ws.Cells.Item(531, 7).Value = "Betarraga"
$ws.Cells.Item(531, 8).Value = "Sin especificar"
$ws.Cells.Item(531, 9).Value = "Primera"
$ws.Cells.Item(531, 10).Value = 60
$ws.Cells.Item(531, 11).Value = 9000
$ws.Cells.Item(531, 12).Value = 9000
$ws.Cells.Item(531, 13).Value = 9000
$ws.Cells.Item(531, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(531, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(531, 16).Value = 750
$ws.Cells.Item(531, 17).Value = 12
$ws.Cells.Item(531, 18).Value = "Hortaliza"
